$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.242.23'
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").Value = '1.605.40'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.37'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("E6").Value = '  +0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.484'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.42%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.25%  '

$ws.Range("E9").Value = '  -0.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.21'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("E11").Value = '  +1.52%  '

$ws.Range("D12").Value = '1.828.30'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").Value = '1.600.39'
$ws.Range("E13").Value = '  -1.15%  '

$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.516'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.57%  '

$ws.Range("D16").Value = '26.263.93'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.33'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.78%  '

$ws.Range("E18").Value = '  -0.17%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.95'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.90%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.28'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.27'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.02'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("E24").Value = '  +11.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.88'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.67%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.122'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.21'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.55'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.86%  '

$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("E33").Value = '  -3.85%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.43'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.49%  '

$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.72%  '

$ws.Range("D36").Value = '1.146.66'
$ws.Range("E36").Value = '  +3.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0162'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +5.82%  '

$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.33'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.786'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.63%  '

$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.783'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.44%  '

$ws.Range("E43").Value = '  +0.81%  '

$ws.Range("D44").Value = '1.742.97'
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.04'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.16%  '

$ws.Range("E46").Value = '  -3.94%  '

$ws.Range("E47").Value = '  +0.73%  '

$ws.Range("E48").Value = '  -0.54%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₇0975'
$ws.Range("E49").Value = '  -7.74%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.408'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.22%  '
